$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Row 669/670: CTO earned 1.25 each (VL(3-0-0) / blank particulars rows) ---
$ws.Cells.Item(669, 3).Value = 1.25
$ws.Cells.Item(670, 3).Value = 1.25

# --- 2. Insert a new row at sheet row 671 (new "2024" year-marker row), pushing
#        everything below down by one (old row 671 -> new row 672, etc.) ---
$ws.Rows.Item(671).Insert()

# New row 671 should look like a normal (blank) data row, formatted like the
# row right below it (old row 671, now at 672) for columns B-K, and like the
# other year-header cells (e.g. A651 = "2023") for column A.
for ($c = 1; $c -le 11; $c++) {
    $src = $ws.Cells.Item(672, $c)
    $dst = $ws.Cells.Item(671, $c)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# A671 = "2024" marker (as text, matching the other year-header cells)
$ws.Cells.Item(671, 1).NumberFormat = "@"
$ws.Cells.Item(671, 1).Value = "2024"
$ws.Cells.Item(651, 1).Copy()
$ws.Cells.Item(671, 1).PasteSpecial(-4122)  # xlPasteFormats (brings quotePrefix style)

# G671 keeps the table's calculated-column formula
$ws.Cells.Item(671, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 3. New row 672 (old row 671, which held 1/1/2024 with no particulars) now
#        also records an SL(1-0-0) CTO usage, same pattern as other SL rows ---
$ws.Cells.Item(672, 2).Value = "SL(1-0-0)"
$ws.Cells.Item(672, 8).Value = 1
$ws.Cells.Item(672, 11).Value = 45300
$ws.Cells.Item(668, 11).Copy()
$ws.Cells.Item(672, 11).PasteSpecial(-4122)  # xlPasteFormats (date style)

# --- 4. Resize Table1 to include the newly-inserted row ---
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K832"))

# Resize can (re)write the totals/last-row calculated-column formula using an
# "[@EARNED]" structured reference that this engine fails to evaluate; put the
# canonical formula back so the cached value recalculates correctly.
$ws.Cells.Item(832, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 5. Restore the cursor/selection state recorded in the workbook ---
$ws.Activate()
$ws.Range("E667").Select()

$wb.Application.Calculate()
